$d = $word.ActiveDocument

$d.Content.Find.Execute("10+48=", $true, $true, $false, $false, $false, $true, 1, $false, "78-38=", 2) | Out-Null
$d.Content.Find.Execute("7+87=", $true, $true, $false, $false, $false, $true, 1, $false, "61-59=", 2) | Out-Null
$d.Content.Find.Execute("18+16=", $true, $true, $false, $false, $false, $true, 1, $false, "91-39=", 2) | Out-Null
$d.Content.Find.Execute("54-47=", $true, $true, $false, $false, $false, $true, 1, $false, "26+73=", 2) | Out-Null
$d.Content.Find.Execute("17-15=", $true, $true, $false, $false, $false, $true, 1, $false, "68+27=", 2) | Out-Null
$d.Content.Find.Execute("53-2=", $true, $true, $false, $false, $false, $true, 1, $false, "45+38=", 2) | Out-Null
$d.Content.Find.Execute("55-37=", $true, $true, $false, $false, $false, $true, 1, $false, "12+4=", 2) | Out-Null
$d.Content.Find.Execute("46-37=", $true, $true, $false, $false, $false, $true, 1, $false, "27+47=", 2) | Out-Null
$d.Content.Find.Execute("63+36=", $true, $true, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("75-75=", $true, $true, $false, $false, $false, $true, 1, $false, "45-4=", 2) | Out-Null
$d.Content.Find.Execute("46-0=", $true, $true, $false, $false, $false, $true, 1, $false, "8-6=", 2) | Out-Null
$d.Content.Find.Execute("7+11=", $true, $true, $false, $false, $false, $true, 1, $false, "87-25=", 2) | Out-Null
$d.Content.Find.Execute("86-80=", $true, $true, $false, $false, $false, $true, 1, $false, "62-25=", 2) | Out-Null
$d.Content.Find.Execute("3+1=", $true, $true, $false, $false, $false, $true, 1, $false, "93-13=", 2) | Out-Null
$d.Content.Find.Execute("55-11=", $true, $true, $false, $false, $false, $true, 1, $false, "6+13=", 2) | Out-Null
$d.Content.Find.Execute("0+64=", $true, $true, $false, $false, $false, $true, 1, $false, "31+29=", 2) | Out-Null
$d.Content.Find.Execute("87-75=", $true, $true, $false, $false, $false, $true, 1, $false, "25+39=", 2) | Out-Null
$d.Content.Find.Execute("99-19=", $true, $true, $false, $false, $false, $true, 1, $false, "49-34=", 2) | Out-Null
$d.Content.Find.Execute("10+43=", $true, $true, $false, $false, $false, $true, 1, $false, "34+16=", 2) | Out-Null
$d.Content.Find.Execute("60+3=", $true, $true, $false, $false, $false, $true, 1, $false, "4+95=", 2) | Out-Null
$d.Content.Find.Execute("82-48=", $true, $true, $false, $false, $false, $true, 1, $false, "30+50=", 2) | Out-Null
$d.Content.Find.Execute("28+52=", $true, $true, $false, $false, $false, $true, 1, $false, "41+7=", 2) | Out-Null
$d.Content.Find.Execute("90-17=", $true, $true, $false, $false, $false, $true, 1, $false, "27-23=", 2) | Out-Null
$d.Content.Find.Execute("16+20=", $true, $true, $false, $false, $false, $true, 1, $false, "31+54=", 2) | Out-Null
$d.Content.Find.Execute("52-32=", $true, $true, $false, $false, $false, $true, 1, $false, "60-38=", 2) | Out-Null
$d.Content.Find.Execute("28-13=", $true, $true, $false, $false, $false, $true, 1, $false, "0+63=", 2) | Out-Null
$d.Content.Find.Execute("36+37=", $true, $true, $false, $false, $false, $true, 1, $false, "70+10=", 2) | Out-Null
$d.Content.Find.Execute("74-6=", $true, $true, $false, $false, $false, $true, 1, $false, "18+28=", 2) | Out-Null
$d.Content.Find.Execute("84-52=", $true, $true, $false, $false, $false, $true, 1, $false, "99-13=", 2) | Out-Null
$d.Content.Find.Execute("39+26=", $true, $true, $false, $false, $false, $true, 1, $false, "97-80=", 2) | Out-Null
$d.Content.Find.Execute("44+37=", $true, $true, $false, $false, $false, $true, 1, $false, "33+3=", 2) | Out-Null
$d.Content.Find.Execute("26+59=", $true, $true, $false, $false, $false, $true, 1, $false, "67+15=", 2) | Out-Null
$d.Content.Find.Execute("71-68=", $true, $true, $false, $false, $false, $true, 1, $false, "79-67=", 2) | Out-Null
$d.Content.Find.Execute("73+4=", $true, $true, $false, $false, $false, $true, 1, $false, "59-3=", 2) | Out-Null
$d.Content.Find.Execute("37+18=", $true, $true, $false, $false, $false, $true, 1, $false, "28+64=", 2) | Out-Null
$d.Content.Find.Execute("25+6=", $true, $true, $false, $false, $false, $true, 1, $false, "8+77=", 2) | Out-Null
$d.Content.Find.Execute("46-35=", $true, $true, $false, $false, $false, $true, 1, $false, "94-36=", 2) | Out-Null
$d.Content.Find.Execute("86+6=", $true, $true, $false, $false, $false, $true, 1, $false, "32-8=", 2) | Out-Null
$d.Content.Find.Execute("48-6=", $true, $true, $false, $false, $false, $true, 1, $false, "53-14=", 2) | Out-Null
$d.Content.Find.Execute("84-81=", $true, $true, $false, $false, $false, $true, 1, $false, "37+43=", 2) | Out-Null
$d.Content.Find.Execute("96-41=", $true, $true, $false, $false, $false, $true, 1, $false, "5+86=", 2) | Out-Null
$d.Content.Find.Execute("87-67=", $true, $true, $false, $false, $false, $true, 1, $false, "12+78=", 2) | Out-Null
$d.Content.Find.Execute("56-26=", $true, $true, $false, $false, $false, $true, 1, $false, "33+36=", 2) | Out-Null
$d.Content.Find.Execute("2+67=", $true, $true, $false, $false, $false, $true, 1, $false, "37+48=", 2) | Out-Null
$d.Content.Find.Execute("33+20=", $true, $true, $false, $false, $false, $true, 1, $false, "2+61=", 2) | Out-Null
$d.Content.Find.Execute("7+44=", $true, $true, $false, $false, $false, $true, 1, $false, "70+26=", 2) | Out-Null
$d.Content.Find.Execute("10+25=", $true, $true, $false, $false, $false, $true, 1, $false, "32+51=", 2) | Out-Null
$d.Content.Find.Execute("94-30=", $true, $true, $false, $false, $false, $true, 1, $false, "66+22=", 2) | Out-Null
$d.Content.Find.Execute("59-42=", $true, $true, $false, $false, $false, $true, 1, $false, "45-44=", 2) | Out-Null
$d.Content.Find.Execute("93-33=", $true, $true, $false, $false, $false, $true, 1, $false, "50-42=", 2) | Out-Null
$d.Content.Find.Execute("85+9=", $true, $true, $false, $false, $false, $true, 1, $false, "0+44=", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $true, $false, $false, $false, $true, 1, $false, "21-15=", 2) | Out-Null
$d.Content.Find.Execute("54+12=", $true, $true, $false, $false, $false, $true, 1, $false, "7+18=", 2) | Out-Null
$d.Content.Find.Execute("4+93=", $true, $true, $false, $false, $false, $true, 1, $false, "36+21=", 2) | Out-Null
$d.Content.Find.Execute("60+32=", $true, $true, $false, $false, $false, $true, 1, $false, "85-58=", 2) | Out-Null
$d.Content.Find.Execute("64+12=", $true, $true, $false, $false, $false, $true, 1, $false, "12+0=", 2) | Out-Null
$d.Content.Find.Execute("85-34=", $true, $true, $false, $false, $false, $true, 1, $false, "67-39=", 2) | Out-Null
$d.Content.Find.Execute("44+19=", $true, $true, $false, $false, $false, $true, 1, $false, "4+39=", 2) | Out-Null
$d.Content.Find.Execute("67+5=", $true, $true, $false, $false, $false, $true, 1, $false, "46-18=", 2) | Out-Null
$d.Content.Find.Execute("73+22=", $true, $true, $false, $false, $false, $true, 1, $false, "72-32=", 2) | Out-Null
$d.Content.Find.Execute("1+65=", $true, $true, $false, $false, $false, $true, 1, $false, "95-61=", 2) | Out-Null
$d.Content.Find.Execute("3-0=", $true, $true, $false, $false, $false, $true, 1, $false, "50-46=", 2) | Out-Null
$d.Content.Find.Execute("21+24=", $true, $true, $false, $false, $false, $true, 1, $false, "71-46=", 2) | Out-Null
$d.Content.Find.Execute("39+0=", $true, $true, $false, $false, $false, $true, 1, $false, "97-76=", 2) | Out-Null
$d.Content.Find.Execute("32-3=", $true, $true, $false, $false, $false, $true, 1, $false, "7+33=", 2) | Out-Null
$d.Content.Find.Execute("93-78=", $true, $true, $false, $false, $false, $true, 1, $false, "5+88=", 2) | Out-Null
$d.Content.Find.Execute("64-18=", $true, $true, $false, $false, $false, $true, 1, $false, "89+1=", 2) | Out-Null
$d.Content.Find.Execute("92-35=", $true, $true, $false, $false, $false, $true, 1, $false, "51+32=", 2) | Out-Null
$d.Content.Find.Execute("90-14=", $true, $true, $false, $false, $false, $true, 1, $false, "12+35=", 2) | Out-Null
$d.Content.Find.Execute("22+47=", $true, $true, $false, $false, $false, $true, 1, $false, "79-13=", 2) | Out-Null
$d.Content.Find.Execute("10+46=", $true, $true, $false, $false, $false, $true, 1, $false, "68-33=", 2) | Out-Null
$d.Content.Find.Execute("55-7=", $true, $true, $false, $false, $false, $true, 1, $false, "15+21=", 2) | Out-Null
$d.Content.Find.Execute("49+15=", $true, $true, $false, $false, $false, $true, 1, $false, "57-28=", 2) | Out-Null
$d.Content.Find.Execute("40+32=", $true, $true, $false, $false, $false, $true, 1, $false, "57-50=", 2) | Out-Null
$d.Content.Find.Execute("86-64=", $true, $true, $false, $false, $false, $true, 1, $false, "1+47=", 2) | Out-Null
$d.Content.Find.Execute("77-70=", $true, $true, $false, $false, $false, $true, 1, $false, "70-22=", 2) | Out-Null
$d.Content.Find.Execute("82-10=", $true, $true, $false, $false, $false, $true, 1, $false, "15+57=", 2) | Out-Null
$d.Content.Find.Execute("9+5=", $true, $true, $false, $false, $false, $true, 1, $false, "23+0=", 2) | Out-Null
$d.Content.Find.Execute("45-18=", $true, $true, $false, $false, $false, $true, 1, $false, "12+7=", 2) | Out-Null
$d.Content.Find.Execute("69-38=", $true, $true, $false, $false, $false, $true, 1, $false, "17-1=", 2) | Out-Null
$d.Content.Find.Execute("98-13=", $true, $true, $false, $false, $false, $true, 1, $false, "73-6=", 2) | Out-Null
$d.Content.Find.Execute("69+5=", $true, $true, $false, $false, $false, $true, 1, $false, "4+36=", 2) | Out-Null
$d.Content.Find.Execute("7+59=", $true, $true, $false, $false, $false, $true, 1, $false, "66-43=", 2) | Out-Null
$d.Content.Find.Execute("43-23=", $true, $true, $false, $false, $false, $true, 1, $false, "81+3=", 2) | Out-Null
$d.Content.Find.Execute("60-25=", $true, $true, $false, $false, $false, $true, 1, $false, "98-49=", 2) | Out-Null
$d.Content.Find.Execute("53-20=", $true, $true, $false, $false, $false, $true, 1, $false, "27+62=", 2) | Out-Null
$d.Content.Find.Execute("28+7=", $true, $true, $false, $false, $false, $true, 1, $false, "64-60=", 2) | Out-Null
$d.Content.Find.Execute("51+26=", $true, $true, $false, $false, $false, $true, 1, $false, "24+39=", 2) | Out-Null
$d.Content.Find.Execute("44+3=", $true, $true, $false, $false, $false, $true, 1, $false, "80+0=", 2) | Out-Null
$d.Content.Find.Execute("95-8=", $true, $true, $false, $false, $false, $true, 1, $false, "84-58=", 2) | Out-Null
$d.Content.Find.Execute("6+92=", $true, $true, $false, $false, $false, $true, 1, $false, "80-5=", 2) | Out-Null
$d.Content.Find.Execute("70-44=", $true, $true, $false, $false, $false, $true, 1, $false, "36-32=", 2) | Out-Null
$d.Content.Find.Execute("72+19=", $true, $true, $false, $false, $false, $true, 1, $false, "56+16=", 2) | Out-Null
$d.Content.Find.Execute("31+61=", $true, $true, $false, $false, $false, $true, 1, $false, "30+8=", 2) | Out-Null
$d.Content.Find.Execute("32+9=", $true, $true, $false, $false, $false, $true, 1, $false, "45-3=", 2) | Out-Null
$d.Content.Find.Execute("26-19=", $true, $true, $false, $false, $false, $true, 1, $false, "31+52=", 2) | Out-Null
$d.Content.Find.Execute("93-45=", $true, $true, $false, $false, $false, $true, 1, $false, "9+24=", 2) | Out-Null
$d.Content.Find.Execute("19-11=", $true, $true, $false, $false, $false, $true, 1, $false, "20+33=", 2) | Out-Null
$d.Content.Find.Execute("62-4=", $true, $true, $false, $false, $false, $true, 1, $false, "32+53=", 2) | Out-Null
$d.Content.Find.Execute("45+8=", $true, $true, $false, $false, $false, $true, 1, $false, "56+14=", 2) | Out-Null
